# Update "想去人数" (F column) figures across sheets, matching the
# gh-pages data refresh generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 516
$ws1.Range("F5").Value  = 494
$ws1.Range("F6").Value  = 927
$ws1.Range("F7").Value  = 165
$ws1.Range("F8").Value  = 958
$ws1.Range("F9").Value  = 746
$ws1.Range("F10").Value = 196
$ws1.Range("F11").Value = 50
$ws1.Range("F13").Value = 785
$ws1.Range("F14").Value = 256
$ws1.Range("F15").Value = 555
$ws1.Range("F17").Value = 1301
$ws1.Range("F19").Value = 428
$ws1.Range("F20").Value = 1113
$ws1.Range("F21").Value = 2801
$ws1.Range("F22").Value = 1312
$ws1.Range("F23").Value = 654
$ws1.Range("F24").Value = 167
$ws1.Range("F27").Value = 973
$ws1.Range("F28").Value = 319
$ws1.Range("F29").Value = 1146
$ws1.Range("F30").Value = 27
$ws1.Range("F32").Value = 1332

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 58

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 516
$ws4.Range("F7").Value  = 494
$ws4.Range("F13").Value = 927
$ws4.Range("F14").Value = 165
$ws4.Range("F16").Value = 958
$ws4.Range("F17").Value = 746
$ws4.Range("F18").Value = 196
$ws4.Range("F19").Value = 58
$ws4.Range("F20").Value = 50
$ws4.Range("F26").Value = 785
$ws4.Range("F27").Value = 256
$ws4.Range("F28").Value = 555
$ws4.Range("F30").Value = 1301
$ws4.Range("F32").Value = 428
$ws4.Range("F33").Value = 1113
$ws4.Range("F34").Value = 2801
$ws4.Range("F35").Value = 1312
$ws4.Range("F36").Value = 654
$ws4.Range("F37").Value = 167
$ws4.Range("F42").Value = 973
$ws4.Range("F43").Value = 319
$ws4.Range("F44").Value = 1146
$ws4.Range("F45").Value = 27
$ws4.Range("F47").Value = 1332
